$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Revert row 6 (onshore wind) shareweights from 3 to 1 across B:AF
$ws.Range("B6:AF6").Value = 1

# Revert row 7 (solar pv) shareweights from 2 to 1 across B:AF
$ws.Range("B7:AF7").Value = 1

# Update the active selection on the ETS sheet to match the saved view
$ws.Activate()
$ws.Range("B6:AF7").Select()
